$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.389.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "'2.625.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'111.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'325.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").Value = "'19.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "'7.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "'3.038.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "'2.632.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "'49.357.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'267.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").Value = "'68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'25.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "'34.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.77%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "'5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "'129.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("D42").Value = "'22.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "'0.0333"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.63%  "
$ws.Range("D45").Value = "'2.058.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "'3.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  +7.66%  "
$ws.Range("E48").Value = "  -5.83%  "
$ws.Range("D49").Value = "'8.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").Value = "'58.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.45%  "
